$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Level_Kode"
$ws.Range("B1").Value = "Level_Nama"
$ws.Range("C1").ClearContents()
$ws.Range("D1").ClearContents()

$ws.Range("A2").Value = ""
$ws.Range("B2").Value = ""

$ws.Columns.Item(1).ColumnWidth = 11.7109375

$ws.Range("A2:B2").Select()
